# [Map] npc object info 시트 추가
$wb = $excel.ActiveWorkbook

$npcInfo = $wb.Worksheets.Item("MapNpcInfo")

# Capture the "PowderShop mirror" object row (row 4) before it gets removed
# from MapNpcInfo, and the updated CompassUIPosition value for row 3.
$objNpcId     = $npcInfo.Cells.Item(4,1).Value2
$objResKey    = $npcInfo.Cells.Item(4,3).Value2
$objSprite    = $npcInfo.Cells.Item(4,4).Value2
$objCollider  = $npcInfo.Cells.Item(4,5).Value2
$objCompass   = $npcInfo.Cells.Item(4,6).Value2

# Create the new "MapObjectInfo" sheet right after "MapNpcInfo".
$objInfo = $wb.Worksheets.Add([System.Type]::Missing, $npcInfo)
$objInfo.Name = "MapObjectInfo"

# Header rows (type row + field-name row), mirroring MapNpcInfo's columns
# A/C/D/E/F (the CharacterId column B is dropped for this sheet).
$objInfo.Cells.Item(1,1).Value = $npcInfo.Cells.Item(1,1).Value2
$objInfo.Cells.Item(1,2).Value = $npcInfo.Cells.Item(1,3).Value2
$objInfo.Cells.Item(1,3).Value = $npcInfo.Cells.Item(1,4).Value2
$objInfo.Cells.Item(1,4).Value = $npcInfo.Cells.Item(1,5).Value2
$objInfo.Cells.Item(1,5).Value = $npcInfo.Cells.Item(1,6).Value2

$objInfo.Cells.Item(2,1).Value = $npcInfo.Cells.Item(2,1).Value2
$objInfo.Cells.Item(2,2).Value = $npcInfo.Cells.Item(2,3).Value2
$objInfo.Cells.Item(2,3).Value = $npcInfo.Cells.Item(2,4).Value2
$objInfo.Cells.Item(2,4).Value = $npcInfo.Cells.Item(2,5).Value2
$objInfo.Cells.Item(2,5).Value = $npcInfo.Cells.Item(2,6).Value2

# Data row moved from MapNpcInfo.
$objInfo.Cells.Item(3,1).Value = $objNpcId
$objInfo.Cells.Item(3,2).Value = $objResKey
$objInfo.Cells.Item(3,3).Value = $objSprite
$objInfo.Cells.Item(3,4).Value = $objCollider
$objInfo.Cells.Item(3,5).Value = $objCompass

# Remove the now-relocated row from MapNpcInfo.
$npcInfo.Rows.Item(4).Delete()

# Update the CompassUIPosition value for the remaining NpcId row (was 0,220).
$npcInfo.Cells.Item(3,6).Value = "0,145"

# Selections per the authored edit.
$npcInfo.Range("C7").Select()
$objInfo.Range("C5").Select()

$posSheet = $wb.Worksheets.Item("MapNpcPosition")
$posSheet.Range("A4:E4").Select()

$menuSheet = $wb.Worksheets.Item("MapNpcMenu")
$menuSheet.Range("H15").Select()

# MapObjectInfo is the active tab when the workbook is saved.
$objInfo.Activate()
